$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment"); C,D,E,F shift left to B,C,D,E.
$ws.Range("B:B").Delete()

# Update the header row text (now shifted into B1:E1).
$ws.Range("B1").Value = "All.jamais.jamais"
$ws.Range("C1").Value = "Males.jamais.jamais"
$ws.Range("D1").Value = "Females.jamais.jamais"
$ws.Range("E1").Value = "Not known / missing.jamais.jamais"
